# Updates the cryptos price/volume table (and swaps the ARBITRUM / Stellar rows)
# to match the latest scrape, per the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.402.14'
$ws.Range('E2').Value = '  +0.13%  '
# Row 3
$ws.Range('D3').Value = '2.105.72'
$ws.Range('E3').Value = '  +4.57%  '
# Row 4
$ws.Range('D4').Value = '''0.9981'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.38%  '
# Row 5
$ws.Range('D5').Value = '''329.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.30%  '
# Row 6
$ws.Range('D6').Value = '''0.9959'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.55%  '
# Row 7
$ws.Range('D7').Value = '''0.5234'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.53%  '
# Row 8
$ws.Range('D8').Value = '''0.4355'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.93%  '
# Row 9
$ws.Range('D9').Value = '''0.08852'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.71%  '
# Row 10
$ws.Range('D10').Value = '''46.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.32%  '
# Row 11
$ws.Range('D11').Value = '''1.163'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.63%  '
# Row 12
$ws.Range('D12').Value = '''24.62'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.58%  '
# Row 13
$ws.Range('D13').Value = '2.089.34'
$ws.Range('E13').Value = '  +3.90%  '
# Row 14
$ws.Range('D14').Value = '''6.740'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.36%  '
# Row 15
$ws.Range('D15').Value = '''7.779'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.23%  '
# Row 16
$ws.Range('D16').Value = '''96.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.20%  '
# Row 17
$ws.Range('D17').Value = '''0.9998'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.23%  '
# Row 18
$ws.Range('D18').Value = '''0.00001130'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.52%  '
# Row 19
$ws.Range('D19').Value = '''0.06628'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.47%  '
# Row 20
$ws.Range('D20').Value = '''18.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.40%  '
# Row 21
$ws.Range('D21').Value = '''0.9997'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.18%  '
# Row 22
$ws.Range('D22').Value = '''6.356'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.41%  '
# Row 23
$ws.Range('D23').Value = '30.454.04'
$ws.Range('E23').Value = '  +0.10%  '
# Row 24
$ws.Range('D24').Value = '''12.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.80%  '
# Row 25
$ws.Range('D25').Value = '''2.324'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.62%  '
# Row 26
$ws.Range('D26').Value = '2.334.08'
$ws.Range('E26').Value = '  +3.83%  '
# Row 27
$ws.Range('D27').Value = '''22.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.34%  '
# Row 28
$ws.Range('D28').Value = '''2.608'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.96%  '
# Row 29
$ws.Range('D29').Value = '''161.61'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.74%  '
# Row 30
$ws.Range('D30').Value = '''132.31'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.68%  '
# Row 31
$ws.Range('D31').Value = '''1.210'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.40%  '
# Row 32
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '''1.708'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +25.36%  '
# Row 33
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').Value = '''0.1072'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.71%  '
# Row 34
$ws.Range('D34').Value = '''6.197'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.03%  '
# Row 35
$ws.Range('D35').Value = '''3.900'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.78%  '
# Row 36
$ws.Range('D36').Value = '''10.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.50%  '
# Row 37
$ws.Range('D37').Value = '''0.02587'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.24%  '
# Row 38
$ws.Range('D38').Value = '''0.06713'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.05%  '
# Row 39
$ws.Range('D39').Value = '''5.491'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.36%  '
# Row 40
$ws.Range('D40').Value = '''12.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.05%  '
# Row 41
$ws.Range('D41').Value = '''0.2267'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.47%  '
# Row 42
$ws.Range('D42').Value = '''0.6829'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.69%  '
# Row 43
$ws.Range('D43').Value = '''1.253'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.78%  '
# Row 44
$ws.Range('D44').Value = '''0.9975'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.35%  '
# Row 45
$ws.Range('D45').Value = '''14.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.65%  '
# Row 46
$ws.Range('D46').Value = '''0.6385'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.59%  '
# Row 47
$ws.Range('D47').Value = '''2.214'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.33%  '
# Row 48
$ws.Range('D48').Value = '''3.616'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.21%  '
# Row 49
$ws.Range('D49').Value = '''1.251'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.04%  '
# Row 50
$ws.Range('D50').Value = '''1.197'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.70%  '
# Row 51
$ws.Range('D51').Value = '''82.08'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.51%  '
